# Applies the "adding new progress as of date 04 nov 2025" update:
#  - Sheet "Training Dashboard": rows 3-23, column H (PERIOD TO EXPIRE) decreases by 1,
#    column I (LAST UPDATE) changes from 03-Nov-2025 to 04-Nov-2025.
#  - Sheet "Exam Dashboard": a new exam result row is inserted before the
#    TOTAL AVERAGE row, and the total average is recalculated.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # Training Dashboard
$ws2 = $wb.Worksheets.Item(2)  # Exam Dashboard

# Helper: write a literal text value into a cell without Excel's automatic
# type conversion (e.g. turning "04-Nov-2025" into a date serial, or
# "81.13%" into a percentage number), while preserving the target cell's
# existing style/format (border, alignment, etc.).
function Set-TextValue($ws, $range, [string]$text) {
    $helper = $ws.Range("Z100")
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $helper.Clear()
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# Training Dashboard: update PERIOD TO EXPIRE (H) and LAST UPDATE (I)
# ---------------------------------------------------------------------
for ($r = 3; $r -le 23; $r++) {
    $hCell = $ws1.Cells.Item($r, 8)
    $hCell.Value = $hCell.Value2 - 1

    $iCell = $ws1.Cells.Item($r, 9)
    Set-TextValue $ws1 $iCell "04-Nov-2025"
}

# ---------------------------------------------------------------------
# Exam Dashboard: widen EXAM column, add new exam row, update total
# ---------------------------------------------------------------------
$ws2.Columns.Item(2).ColumnWidth = 23.17   # renders as stored width 24

# Copy formatting of the current TOTAL AVERAGE row (row 5) down to the
# new row 6 before it has any content, so the shifted row keeps the same
# cell style as the rest of the table.
$ws2.Range("A5:G5").Copy()
$ws2.Range("A6:G6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Move the TOTAL AVERAGE label down to row 6 with the recalculated value.
Set-TextValue $ws2 $ws2.Range("C6") "TOTAL AVERAGE"
Set-TextValue $ws2 $ws2.Range("D6") "79.46%"

# Fill in the new exam entry on row 5.
$ws2.Range("A5").Value = 3
Set-TextValue $ws2 $ws2.Range("B5") "Consignment Shuttle Tv"
Set-TextValue $ws2 $ws2.Range("C5") "29-Oct-2025"
Set-TextValue $ws2 $ws2.Range("D5") "81.13%"
Set-TextValue $ws2 $ws2.Range("E5") "VALID"
Set-TextValue $ws2 $ws2.Range("F5") "Approved Score. date is valid"
